# Generate Report for Handoff
# This script re-applies the localization-status report refresh:
#  - the "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md" entry moves from row 2 to row 4
#    (it is now "Ready for handoff" with a fresh handoff timestamp and an
#    out-of-date handback warning), while "ffff5a98fd62...md" and
#    "ffffff3b1be90b...md" shift up into rows 2 and 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md"
$ws1.Range("B2").Value = "e2e\ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md"
$ws1.Range("G2").Value = "2016-08-16 18:57:20"

$ws1.Range("A3").Value = "ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md"
$ws1.Range("B3").Value = "e2e\ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md"

$ws1.Range("A4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md"
$ws1.Range("B4").Value = "e2e\043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md"
$ws1.Range("E4").Value = "Ready for handoff"
$ws1.Range("F4").Value = "Ready for handoff"
$ws1.Range("G4").Value = "2016-08-16 19:00:26"

# Hyperlinks keep pointing at the same target URLs as before, just the
# row/display text they are attached to changes.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55a18638d58d3fb31535528db9cd6d5f32523dbd/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md", "", "", "e2e\ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a454bda4b5ecf50d36f125308d4ac209d299215/e2e/ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md", "", "", "e2e\ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55a18638d58d3fb31535528db9cd6d5f32523dbd/e2e/ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md", "", "", "e2e\043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md"
$ws2.Range("G2").Value = "309396e2-187e-4ce5-9a72-944cbbf11640.1a2cd9f118dcb939d6178d74844f5128593e0c36.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-16 18:57:15"
$ws2.Range("I2").Value = "309396e2-187e-4ce5-9a72-944cbbf11640.md"
$ws2.Range("J2").Value = "309396e2-187e-4ce5-9a72-944cbbf11640.1a2cd9f118dcb939d6178d74844f5128593e0c36.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-16 18:57:32"

$ws2.Range("A3").Value = "ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md"
$ws2.Range("F3").Value = "True"

$ws2.Range("A4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("F4").Value = "False"
$ws2.Range("G4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.72197f1bdc1eeae94b2b0865ae734c409ee96697.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-08-16 19:00:03"
$ws2.Range("I4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md"
$ws2.Range("J4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.72197f1bdc1eeae94b2b0865ae734c409ee96697.zh-cn.xlf"
$ws2.Range("K4").Value = "2016-08-16 18:59:36"
$ws2.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55a18638d58d3fb31535528db9cd6d5f32523dbd/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15affc1f1d6f265c434e84a2f8c53759fe99255d/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md."

$ws2.Columns.Item(16).ColumnWidth = 39.166666666666664

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55a18638d58d3fb31535528db9cd6d5f32523dbd/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md", "", "", "ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c81827150db2267f6b34cecc5716e27d9b62f7eb/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md", "", "", "309396e2-187e-4ce5-9a72-944cbbf11640.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a454bda4b5ecf50d36f125308d4ac209d299215/e2e/ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md", "", "", "ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/80a861cb3c6a25c55fe2d66798cf71e861a47ba6/e2e/309396e2-187e-4ce5-9a72-944cbbf11640.md", "", "", "309396e2-187e-4ce5-9a72-944cbbf11640.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55a18638d58d3fb31535528db9cd6d5f32523dbd/e2e/ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md", "", "", "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/80a861cb3c6a25c55fe2d66798cf71e861a47ba6/e2e/309396e2-187e-4ce5-9a72-944cbbf11640.md", "", "", "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md"
$ws3.Range("G2").Value = "309396e2-187e-4ce5-9a72-944cbbf11640.1a2cd9f118dcb939d6178d74844f5128593e0c36.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-16 18:57:20"
$ws3.Range("I2").Value = "309396e2-187e-4ce5-9a72-944cbbf11640.md"
$ws3.Range("J2").Value = "309396e2-187e-4ce5-9a72-944cbbf11640.1a2cd9f118dcb939d6178d74844f5128593e0c36.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-16 18:57:39"

$ws3.Range("A3").Value = "ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md"
$ws3.Range("F3").Value = "True"

$ws3.Range("A4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("F4").Value = "False"
$ws3.Range("G4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.72197f1bdc1eeae94b2b0865ae734c409ee96697.de-de.xlf"
$ws3.Range("H4").Value = "2016-08-16 19:00:26"
$ws3.Range("I4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md"
$ws3.Range("J4").Value = "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.72197f1bdc1eeae94b2b0865ae734c409ee96697.de-de.xlf"
$ws3.Range("K4").Value = "2016-08-16 18:59:44"
$ws3.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55a18638d58d3fb31535528db9cd6d5f32523dbd/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15affc1f1d6f265c434e84a2f8c53759fe99255d/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md."

$ws3.Columns.Item(16).ColumnWidth = 39.166666666666664

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55a18638d58d3fb31535528db9cd6d5f32523dbd/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md", "", "", "ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/484ac8eadfd7c74c22ef49042b472db5b62cf625/e2e/043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md", "", "", "309396e2-187e-4ce5-9a72-944cbbf11640.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a454bda4b5ecf50d36f125308d4ac209d299215/e2e/ffff5a98fd62-3905-45fb-b1c2-9da87edb1a24.md", "", "", "ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bb5df79ca8e914f73e145d4044fc2f3908669310/e2e/309396e2-187e-4ce5-9a72-944cbbf11640.md", "", "", "309396e2-187e-4ce5-9a72-944cbbf11640.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55a18638d58d3fb31535528db9cd6d5f32523dbd/e2e/ffffff3b1be90b-93fc-4f62-8334-1ff8cc9c5dac.md", "", "", "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bb5df79ca8e914f73e145d4044fc2f3908669310/e2e/309396e2-187e-4ce5-9a72-944cbbf11640.md", "", "", "043baf2f-cdd5-4ee9-b21d-8e4ee0f5d3e8.md") | Out-Null
